$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C held "bad" dates (year 2002) because of an old, broken validation
# source. Point the rows at the corrected 2020 dates (same day offsets as
# column B) now that validation has moved to a different file.
$ws.Range("C2").Value = 44184
$ws.Range("C3").Value = 44185
$ws.Range("C4").Value = 44186
$ws.Range("C5").Value = 44187
$ws.Range("C6").Value = 44188
$ws.Range("C7").Value = 44189
$ws.Range("C8").Value = 44190
$ws.Range("C9").Value = 44191
$ws.Range("C10").Value = 44192

# Header row value becomes an intentionally "broken" text stand-in, matching
# the existing pattern used in B1 ("11/10/20a20").
$ws.Range("C1").Value = "12/18/20a20"

# Selection moves to C1.
$ws.Range("C1").Select()
